# Added 1.1.0 of term
#
# 1. Bump the "Version" metadata value 1.0.0 -> 1.1.0
# 2. Bump the "Date" metadata value to the new publication timestamp
# 3. Re-apply the wrap/top alignment to the already-bordered cells so the
#    workbook's cell formats carry an explicit alignment flag again
#    (mirrors the applyAlignment="true" added to the two borderId="8" xf
#    records in styles.xml).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-assert the existing vertical-top / wrap-text alignment on every cell
# that already carries it, scoped exactly to the populated cells so no new
# style variants are created for previously-untouched (blank) cells.
$meta.Range("A1:B14").WrapText = $true

$incl = $wb.Worksheets.Item("Include from FSIII")
$incl.Range("A1:C1").WrapText = $true
$incl.Range("A2:C2").WrapText = $true
$incl.Range("A3:B4").WrapText = $true
